$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.182.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.71%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.671.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.27%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'592.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'166.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -6.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.670.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.97%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.38%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -5.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.83%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -5.48%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'37.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -6.06%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000239"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.59%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.289.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.682.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.03%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'67.239.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -4.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -6.73%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'16.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.60%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'483.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -7.87%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.714"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'84.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.05%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -7.59%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000139"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.60%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -5.94%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.47%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -6.51%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -3.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -6.95%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.79%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'31.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.11%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.817.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.99%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.616.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.40%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -7.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.40%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.988"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -6.25%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.87%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -7.63%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -5.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'440.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -9.31%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'48.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.72%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -7.63%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -8.40%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'8.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.76%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.01%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'141.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.83%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'39.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -10.86%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.750.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -6.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0345"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.28%  "
$ws.Range("E51").Style = "Normal"
